$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.205702586464461
$ws.Range("B3").Value = 0.1643932307303566
$ws.Range("B4").Value = 0.1301629405423668
$ws.Range("C4").Value = 2
$ws.Range("B5").Value = 0.109029252404279
$ws.Range("C5").Value = 2
$ws.Range("B6").Value = 0.2349255361509649
$ws.Range("B7").Value = 0.1753554597463773
$ws.Range("B8").Value = 0.1039229826610508
$ws.Range("C8").Value = 2
$ws.Range("B9").Value = 0.1127281548649403
$ws.Range("C9").Value = 2
$ws.Range("B10").Value = 0.2657546054193928
$ws.Range("B11").Value = 0.1842050212129164
$ws.Range("B12").Value = 0.07476217823650765
$ws.Range("B13").Value = 0.1275136930872255
$ws.Range("C13").Value = 2
$ws.Range("B14").Value = 0.1724622471301214
$ws.Range("B15").Value = 0.176267729111987
$ws.Range("B16").Value = 0.1949085160080017
$ws.Range("B17").Value = 0.2686615794129055
$ws.Range("B18").Value = 0.09884205815116515
$ws.Range("C18").Value = 2
$ws.Range("B19").Value = 0.09612066273305153
$ws.Range("C19").Value = 2
$ws.Range("B20").Value = 0.2505747425742995
$ws.Range("B21").Value = 0.1882567665545256
$ws.Range("B22").Value = 0.1691909686138251
$ws.Range("B23").Value = 0.1771075965611703
$ws.Range("B24").Value = 0.1777097446305805
$ws.Range("B25").Value = 0.1988330232120193
$ws.Range("B26").Value = 0.1558336284130092
$ws.Range("C26").Value = 2
$ws.Range("B27").Value = 0.151885000149397
$ws.Range("C27").Value = 2
$ws.Range("B28").Value = 0.2646095887525755
$ws.Range("B29").Value = 0.1193517335240628
$ws.Range("C29").Value = 2
$ws.Range("B30").Value = 0.09585285992210585
$ws.Range("C30").Value = 2
$ws.Range("B31").Value = 0.1587746591110122
$ws.Range("C31").Value = 2
$ws.Range("B32").Value = 0.1390706124909216
$ws.Range("C32").Value = 2
$ws.Range("B33").Value = 0.2116139315861498
$ws.Range("B34").Value = 0.2001899161844258
$ws.Range("B35").Value = 0.2959939320284075
$ws.Range("B36").Value = 0.2063818574789965
$ws.Range("B37").Value = 0.181271018807856
$ws.Range("B38").Value = 0.2237414899724128
$ws.Range("B39").Value = 0.2454544098688217
$ws.Range("B40").Value = 0.1473826792247265
$ws.Range("C40").Value = 2
$ws.Range("B41").Value = 0.2968459084409374
$ws.Range("B42").Value = 0.1098208266875632
$ws.Range("C42").Value = 2
$ws.Range("B43").Value = 0.1696215075973009
$ws.Range("B44").Value = 0.1179767179609373
$ws.Range("C44").Value = 2
$ws.Range("B45").Value = 0.106401399922419
$ws.Range("C45").Value = 2
$ws.Range("B46").Value = 0.1978244236789283
$ws.Range("B47").Value = 0.2175725242527419
$ws.Range("B48").Value = 0.1849190056139085
$ws.Range("B49").Value = 0.2045798442404839
$ws.Range("B50").Value = 0.0304819300579829
$ws.Range("C50").Value = 3
$ws.Range("B51").Value = 0.03275573029001756
